$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the Runmode ("Y") column for the existing DRA0010-DRA0013 rows
$ws.Range("D6").Value = "Y"
$ws.Range("D7").Value = "Y"
$ws.Range("D8").Value = "Y"
$ws.Range("D9").Value = "Y"

# Add the new DRA5 test case row (written in this order so new shared
# strings land in the same index order as the authored workbook)
$ws.Range("A10").Value = "DRA5"
$ws.Range("C10").Value = "Verify that ProfileFlyout links"
$ws.Range("B10").Value = "OPQA-TBD"
$ws.Range("D10").Value = "Y"

# Copy formatting from the row above so the new row matches the rest of the table
$ws.Range("A9:E9").Copy()
$ws.Range("A10:E10").PasteSpecial(-4122)

# Update the active selection / scrolled position to match the authored state
$ws.Range("C14").Select()
$excel.ActiveWindow.ScrollRow = 5
